$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("S2c - pyDamage analysis")

$values = @{
    4  = 44
    5  = 24
    6  = 98.8
    7  = 98.09999999999999
    8  = 72.90000000000001
    9  = 99.40000000000001
    10 = 86.3
    11 = 60.4
    12 = 85.3
    13 = 87.59999999999999
    14 = 71.2
    15 = 70.40000000000001
    16 = 16.7
    17 = 90.09999999999999
    18 = 25.8
    19 = 90.40000000000001
    20 = 60
    21 = 86.5
    22 = 66.5
    23 = 91.09999999999999
    24 = 78
    25 = 91
    26 = 49.7
    27 = 97.5
    28 = 95.7
    29 = 97.40000000000001
    30 = 94.90000000000001
    31 = 90.59999999999999
    32 = 62
    33 = 78.09999999999999
    34 = 48.6
    35 = 95.40000000000001
    36 = 88.3
    37 = 88.2
    38 = 86.7
    39 = 89.40000000000001
    40 = 91.8
    41 = 94.09999999999999
    42 = 85.90000000000001
    43 = 73
    44 = 15.2
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 4).Value = $values[$row]
}
